$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "B" and "C" quarter rows within each year group (rows 3/4, 7/8, 11/12, 15/16)
$pairs = @(
    @(3, 4),
    @(7, 8),
    @(11, 12),
    @(15, 16)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $row1 = $ws.Range("A" + $r1 + ":E" + $r1).Value()
    $row2 = $ws.Range("A" + $r2 + ":E" + $r2).Value()

    $ws.Range("A" + $r1 + ":E" + $r1).Value = $row2
    $ws.Range("A" + $r2 + ":E" + $r2).Value = $row1
}

# Delete columns F:G entirely (drop the extra "原盐产销率" / "原盐销售量" columns)
$ws.Range("F1:G17").Delete()
